$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not auto-converted to a number),
# then clear the temporary "@" number-format override so the cell
# keeps the default style it had before (no leftover s="..").
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "70.907.85"
$ws.Range("E2").Value = "  -2.39%  "
Set-TextValue $ws.Range("D3") "3.939.37"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "536.17"
$ws.Range("E5").Value = "  +2.86%  "
Set-TextValue $ws.Range("D6") "147.52"
$ws.Range("E6").Value = "  -0.36%  "
Set-TextValue $ws.Range("D7") "3.931.29"
$ws.Range("E7").Value = "  -2.98%  "
Set-TextValue $ws.Range("D8") "0.684"
$ws.Range("E8").Value = "  -6.20%  "
$ws.Range("E9").Value = "  +0.01%  "
Set-TextValue $ws.Range("D10") "0.738"
$ws.Range("E10").Value = "  -5.50%  "
Set-TextValue $ws.Range("D11") "0.165"
$ws.Range("E11").Value = "  -6.33%  "
Set-TextValue $ws.Range("D12") "54.77"
$ws.Range("E12").Value = "  +12.75%  "
Set-TextValue $ws.Range("D13") "0.0000316"
$ws.Range("E13").Value = "  -4.51%  "
Set-TextValue $ws.Range("D14") "10.58"
$ws.Range("E14").Value = "  -4.77%  "
Set-TextValue $ws.Range("D15") "4.572.09"
$ws.Range("E15").Value = "  -2.89%  "
Set-TextValue $ws.Range("D16") "3.940.08"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("E17").Value = "  -3.63%  "
Set-TextValue $ws.Range("D18") "20.48"
$ws.Range("E18").Value = "  -4.02%  "
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("E20").Value = "  -5.52%  "
Set-TextValue $ws.Range("D21") "70.879.62"
$ws.Range("E21").Value = "  -2.29%  "
Set-TextValue $ws.Range("D22") "419.89"
$ws.Range("E22").Value = "  -6.80%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D23") "3.59"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D24") "97.35"
$ws.Range("E24").Value = "  -7.46%  "
Set-TextValue $ws.Range("D25") "4.24"
$ws.Range("E25").Value = "  +5.51%  "
Set-TextValue $ws.Range("D26") "14.53"
$ws.Range("E26").Value = "  -3.90%  "
Set-TextValue $ws.Range("D27") "11.29"
$ws.Range("E27").Value = "  -0.95%  "
Set-TextValue $ws.Range("D28") "3.85"
$ws.Range("E28").Value = "  +16.51%  "
Set-TextValue $ws.Range("D29") "10.72"
$ws.Range("E29").Value = "  -3.42%  "
Set-TextValue $ws.Range("D30") "5.89"
$ws.Range("E30").Value = "  +0.86%  "
Set-TextValue $ws.Range("D31") "36.40"
$ws.Range("E31").Value = "  -4.62%  "
Set-TextValue $ws.Range("D32") "7.68"
$ws.Range("E32").Value = "  +16.04%  "
Set-TextValue $ws.Range("D33") "50.85"
$ws.Range("E33").Value = "  +20.11%  "
$ws.Range("E34").Value = "  +0.93%  "
Set-TextValue $ws.Range("D35") "13.32"
$ws.Range("E35").Value = "  -2.68%  "
Set-TextValue $ws.Range("D36") "676.67"
$ws.Range("E36").Value = "  -0.20%  "
Set-TextValue $ws.Range("D37") "65.58"
$ws.Range("E37").Value = "  -3.68%  "
Set-TextValue $ws.Range("D38") "0.443"
$ws.Range("E38").Value = "  +2.87%  "
Set-TextValue $ws.Range("D39") "0.0₃0810"
$ws.Range("E39").Value = "  -6.28%  "
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  -2.78%  "
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +0.06%  "
Set-TextValue $ws.Range("D44") "0.0481"
$ws.Range("E44").Value = "  -3.76%  "
Set-TextValue $ws.Range("D45") "3.18"
$ws.Range("E45").Value = "  -1.23%  "
Set-TextValue $ws.Range("D46") "10.13"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  -5.81%  "
Set-TextValue $ws.Range("D48") "2.64"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -2.91%  "
Set-TextValue $ws.Range("D51") "144.17"
$ws.Range("E51").Value = "  -0.23%  "
